$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New texts used in this edit ---
$txtN11 = "Liberacion de equipo de frontera, prueba de comunicacion de sensor de fondo, validacion de funcionamiento de motor con el motor de banco encontrados, validacion del PQM valores anormales de potencia activa , y factor de potencia entregado por la red"
$txtR9  = "Compensatorio"
$txtN17 = "Gestion de Viaje, compensatorio, y dias de la familia "

# --- Row 11 : taller row, new N11 note, O11 changes 30 -> 95, S11:U11 cleared ---
$ws.Rows(11).RowHeight = 96
$ws.Range("N11").Value = $txtN11
$ws.Range("O11").Value = 95
$ws.Range("S11:U11").Value = ""

# --- Row 9 : O9:U9 block gets cleared, R9 gets the "Compensatorio" note ---
$ws.Range("O9:U9").Value = ""
$ws.Range("R9").Value = $txtR9

# --- Row 13 : S13:U13 cleared ---
$ws.Range("S13:U13").Value = ""

# --- Row 17 : N17 gets the "dias de la familia" note (wrap text), O17 changes 30 -> 5 ---
$ws.Range("N17").Value = $txtN17
$ws.Range("N17").WrapText = $true
$ws.Range("O17").Value = 5

# --- Row 18 : N18 alignment now wraps text (style only, stays empty) ---
$ws.Range("N18").WrapText = $true

# --- Column N width grows to fit the new long notes ---
$ws.Columns(14).ColumnWidth = 32.33

# --- Selection / scroll position as left by the author ---
$ws.Range("R11:R12").Select()
$excel.ActiveWindow.ScrollRow = 7
